# Fruta / hortaliza, semanal
# Insert three new weekly price rows for "Black Amber" (Femacal de La Calera - Ciruela)
# right after the existing "Black Amber 10 kilos" block (before the old row 62),
# which pushes all subsequent rows down by 3 (old row 135 -> new row 138).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 62 (existing rows 62..135 shift down to 65..138)
$ws.Rows("62:64").Insert()

# New row 62
$ws.Range("A62").Value = 3
$ws.Range("B62").Value = "Femacal de La Calera"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = "2022-01-18"
$ws.Range("E62").Value = 5
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100103
$ws.Range("H62").Value = "Frutos de hueso (carozo)"
$ws.Range("I62").Value = 100103002
$ws.Range("J62").Value = "Ciruela"
$ws.Range("K62").Value = "Black Amber"
$ws.Range("L62").Value = "Especial"
$ws.Range("M62").Value = 90
$ws.Range("N62").Value = 13000
$ws.Range("O62").Value = 13000
$ws.Range("P62").Value = 13000
$ws.Range("Q62").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R62").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S62").Value = 1300
$ws.Range("T62").Value = 10

# New row 63
$ws.Range("A63").Value = 3
$ws.Range("B63").Value = "Femacal de La Calera"
$ws.Range("C63").Value = "Coquimbo"
$ws.Range("D63").Value = "2022-01-18"
$ws.Range("E63").Value = 5
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100103
$ws.Range("H63").Value = "Frutos de hueso (carozo)"
$ws.Range("I63").Value = 100103002
$ws.Range("J63").Value = "Ciruela"
$ws.Range("K63").Value = "Black Amber"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 90
$ws.Range("N63").Value = 12000
$ws.Range("O63").Value = 12000
$ws.Range("P63").Value = 12000
$ws.Range("Q63").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R63").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S63").Value = 1200
$ws.Range("T63").Value = 10

# New row 64
$ws.Range("A64").Value = 3
$ws.Range("B64").Value = "Femacal de La Calera"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = "2022-01-18"
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100103
$ws.Range("H64").Value = "Frutos de hueso (carozo)"
$ws.Range("I64").Value = 100103002
$ws.Range("J64").Value = "Ciruela"
$ws.Range("K64").Value = "Black Amber"
$ws.Range("L64").Value = "Segunda"
$ws.Range("M64").Value = 90
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 10000
$ws.Range("P64").Value = 10000
$ws.Range("Q64").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R64").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S64").Value = 1000
$ws.Range("T64").Value = 10
